{"js": "// Add a new paragraph after the existing content with the red (FF0000)\n// text \"Check list guia para entrevista\" \u2014 matching the diff that inserts\n// a \"Check list guia para entrevista\" paragraph right after \"F\u00c9RIAS\".\n\nconst body = context.document.body;\n\n// Create a new, empty paragraph at the very end of the document body.\nconst newPara = body.insertParagraph(\"\", Word.InsertLocation.end);\n\n// Color the (still empty) paragraph red first so the paragraph mark\n// itself picks up the red run formatting (mirrors w:pPr/w:rPr/w:color).\nnewPara.font.color = \"#FF0000\";\nawait context.sync();\n\n// Now insert the actual text at the end of that paragraph and make sure\n// the inserted run is red as well (mirrors w:r/w:rPr/w:color).\nconst insertedRange = newPara.getRange(Word.RangeLocation.end);\ninsertedRange.insertText(\"Check list guia para entrevista\", Word.InsertLocation.replace);\ninsertedRange.font.color = \"#FF0000\";\n\nawait context.sync();\n", "ps1": "# Add a new paragraph after the existing content with the red (FF0000)\n# text \"Check list guia para entrevista\" \u2014 matching the diff that inserts\n# a \"Check list guia para entrevista\" paragraph right after \"F\u00c9RIAS\".\n\n$d = $word.ActiveDocument\n\n# Collapse to the end of the document and insert a new paragraph mark there.\n$end = $d.Content\n$end.Collapse(0)\n$end.InsertParagraphAfter()\n\n# The newly created paragraph is now the last paragraph in the document.\n$newPara = $d.Paragraphs.Last\n$newPara.Range.Text = \"Check list guia para entrevista\"\n\n# wdColorRed (255 / 0x0000FF as a BGR COM color value) -> OOXML FF0000.\n$newPara.Range.Font.Color = 255\n"}
